$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row ("adding some more stuff")
$ws.Range("A5").Value = "First Kranthi"
$ws.Range("B5").Value = "Last Kumar"
$ws.Range("C5").Value = "Ph#89"
$ws.Range("D5").Value = "Address123"

# The longer values in column A/B make Excel's "best fit" column width grow;
# nudge both columns so their widths reflect the new, longer content.
$ws.Columns.Item(1).ColumnWidth = 8.65
$ws.Columns.Item(2).ColumnWidth = 8.5
